# Perbaikan Antrian Device Presensi
# Fix attendance status values for a few students and update the summary counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (AHYATUL HUSNA): Hadir -> Sakit
$ws.Range("D4").Value = "Sakit"

# Row 7 (AMA DYA MOZA): Hadir -> Alpha
$ws.Range("D7").Value = "Alpha"

# Row 8 (ANAFAUL RISTA RAMADANI): Hadir -> Alpha
$ws.Range("D8").Value = "Alpha"

# Row 16 (FARIS NUR ARIFIN): Sakit -> Hadir
$ws.Range("D16").Value = "Hadir"

# Row 30 (WIBI HASAN WIRATAMA): Sakit -> Hadir
$ws.Range("D30").Value = "Hadir"

# Update the summary counts at the bottom of the sheet.
$ws.Range("A36").Value = "Hadir: 27"
$ws.Range("A38").Value = "Sakit: 1"
$ws.Range("A39").Value = "Alpha: 3"
